{"js": "// Update the date line and the 25 division problems in the practice-sheet\n// table. Each table cell is addressed by (row, column) rather than by a\n// global text search so that values which collide with each other across\n// cells (e.g. \"31\u00f73=\" is both an old value in one cell and a new value in\n// another) are never mismatched.\n\n// 1) Update the date heading paragraph.\nconst dateResults = context.document.body.search(\"2024-07-13 Saturday\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"2024-07-14 Sunday\", Word.InsertLocation.replace);\n}\n\n// 2) Update the division problems inside the table, cell by cell.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// (row, column, old text, new text) - only the five \"problem\" rows (the\n// other rows in the table are blank answer rows and are left untouched).\nconst updates = [\n  [0, 0, \"63\u00f72=\", \"19\u00f79=\"],\n  [0, 1, \"77\u00f77=\", \"87\u00f75=\"],\n  [0, 2, \"82\u00f79=\", \"31\u00f73=\"],\n  [0, 3, \"28\u00f79=\", \"59\u00f77=\"],\n  [0, 4, \"84\u00f72=\", \"67\u00f79=\"],\n  [4, 0, \"72\u00f77=\", \"95\u00f77=\"],\n  [4, 1, \"34\u00f79=\", \"93\u00f78=\"],\n  [4, 2, \"93\u00f79=\", \"39\u00f79=\"],\n  [4, 3, \"13\u00f76=\", \"47\u00f76=\"],\n  [4, 4, \"49\u00f79=\", \"94\u00f74=\"],\n  [8, 0, \"58\u00f79=\", \"35\u00f72=\"],\n  [8, 1, \"10\u00f77=\", \"77\u00f77=\"],\n  [8, 2, \"68\u00f73=\", \"16\u00f73=\"],\n  [8, 3, \"84\u00f74=\", \"17\u00f78=\"],\n  [8, 4, \"80\u00f73=\", \"23\u00f78=\"],\n  [12, 0, \"83\u00f77=\", \"79\u00f78=\"],\n  [12, 1, \"31\u00f73=\", \"65\u00f74=\"],\n  [12, 2, \"44\u00f74=\", \"20\u00f75=\"],\n  [12, 3, \"37\u00f72=\", \"39\u00f79=\"],\n  [12, 4, \"46\u00f75=\", \"39\u00f76=\"],\n  [16, 0, \"49\u00f77=\", \"84\u00f73=\"],\n  [16, 1, \"54\u00f79=\", \"91\u00f78=\"],\n  [16, 2, \"63\u00f75=\", \"64\u00f76=\"],\n  [16, 3, \"33\u00f76=\", \"56\u00f77=\"],\n  [16, 4, \"58\u00f77=\", \"49\u00f73=\"],\n];\n\nfor (const [row, col, oldText, newText] of updates) {\n  const cell = table.getCell(row, col);\n  const cellResults = cell.body.search(oldText, { matchCase: true });\n  cellResults.load(\"items\");\n  await context.sync();\n  if (cellResults.items.length > 0) {\n    cellResults.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division problems in the practice-sheet\n# table. Each location is addressed structurally (paragraph index / table\n# cell row+column) and updated by assigning Range.Text directly, rather\n# than via a document-wide Find, so that values which collide with each\n# other across cells (e.g. \"31\u00f73=\" is both an old value in one cell and a\n# new value in another) are never mismatched and formatting (fonts, size,\n# alignment) carried by the existing run/paragraph is preserved.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date heading paragraph (first paragraph in the document).\n$d.Paragraphs.Item(1).Range.Text = \"2024-07-14 Sunday\"\n\n# 2) Update the division problems inside the table, cell by cell.\n$tbl = $d.Tables.Item(1)\n\n# (row, column, new text) - 1-based row/column indices; only the five\n# \"problem\" rows are listed (the other rows in the table are blank answer\n# rows and are left untouched).\n$updates = @(\n    @(1, 1, \"19\u00f79=\"),\n    @(1, 2, \"87\u00f75=\"),\n    @(1, 3, \"31\u00f73=\"),\n    @(1, 4, \"59\u00f77=\"),\n    @(1, 5, \"67\u00f79=\"),\n    @(5, 1, \"95\u00f77=\"),\n    @(5, 2, \"93\u00f78=\"),\n    @(5, 3, \"39\u00f79=\"),\n    @(5, 4, \"47\u00f76=\"),\n    @(5, 5, \"94\u00f74=\"),\n    @(9, 1, \"35\u00f72=\"),\n    @(9, 2, \"77\u00f77=\"),\n    @(9, 3, \"16\u00f73=\"),\n    @(9, 4, \"17\u00f78=\"),\n    @(9, 5, \"23\u00f78=\"),\n    @(13, 1, \"79\u00f78=\"),\n    @(13, 2, \"65\u00f74=\"),\n    @(13, 3, \"20\u00f75=\"),\n    @(13, 4, \"39\u00f79=\"),\n    @(13, 5, \"39\u00f76=\"),\n    @(17, 1, \"84\u00f73=\"),\n    @(17, 2, \"91\u00f78=\"),\n    @(17, 3, \"64\u00f76=\"),\n    @(17, 4, \"56\u00f77=\"),\n    @(17, 5, \"49\u00f73=\")\n)\n\nforeach ($u in $updates) {\n    $row = $u[0]\n    $col = $u[1]\n    $newText = $u[2]\n    $tbl.Cell($row, $col).Range.Text = $newText\n}\n"}
